# chore: update Sheets via scheduled runner
# Refresh market-price derived columns (currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# for the rows whose upstream prices changed, across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2429.1724
$ws.Range("I15").Value = 2429.1724
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7287.5172
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -7118.5172
$ws.Range("N69").ClearContents()
$ws.Range("H69").Value = 2999
$ws.Range("I69").Value = 2999
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 8997
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -8123
$ws.Range("N72").ClearContents()
$ws.Range("H72").Value = 2999
$ws.Range("I72").Value = 2999
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 26991
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -22623
$ws.Range("H113").Value = 4477.6665
$ws.Range("I113").Value = 3200
$ws.Range("J113").Value = 5116.5
$ws.Range("K113").Value = 3200
$ws.Range("L113").Value = 5116.5
$ws.Range("M113").Value = 54
$ws.Range("H132").Value = 2312.5
$ws.Range("I132").Value = 2313.889
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 6941.667
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = -4411.667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4449.0557
$ws.Range("I32").Value = 3905.5312
$ws.Range("J32").Value = 8797.25
$ws.Range("K32").Value = 3905.5312
$ws.Range("L32").Value = 8797.25
$ws.Range("M32").Value = -3618.5312
$ws.Range("H74").Value = 1261
$ws.Range("I74").Value = 1198.4667
$ws.Range("J74").Value = 2199
$ws.Range("K74").Value = 1198.4667
$ws.Range("L74").Value = 2199
$ws.Range("M74").Value = -324.4666999999999
$ws.Range("H77").Value = 1261
$ws.Range("I77").Value = 1198.4667
$ws.Range("J77").Value = 2199
$ws.Range("K77").Value = 5992.3335
$ws.Range("L77").Value = 10995
$ws.Range("M77").Value = -1624.3335
$ws.Range("H97").Value = 1543.4445
$ws.Range("I97").Value = 532.8
$ws.Range("J97").Value = 2806.75
$ws.Range("K97").Value = 532.8
$ws.Range("L97").Value = 2806.75
$ws.Range("M97").Value = -36.79999999999995
$ws.Range("H102").Value = 2889.1428
$ws.Range("I102").Value = 2120.6667
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 2120.6667
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -498.6667000000002
$ws.Range("N102").Value = -10744
$ws.Range("H122").Value = 2612.3076
$ws.Range("I122").Value = 2706.25
$ws.Range("J122").Value = 1485
$ws.Range("K122").Value = 8118.75
$ws.Range("L122").Value = 4455
$ws.Range("M122").Value = -5668.75
$ws.Range("H132").Value = 2632.4783
$ws.Range("I132").Value = 2407.45
$ws.Range("J132").Value = 4132.6665
$ws.Range("K132").Value = 7222.349999999999
$ws.Range("L132").Value = 12397.9995
$ws.Range("M132").Value = -4692.349999999999
$ws.Range("N132").Value = -17457.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4249.643
$ws.Range("I20").Value = 4788.8887
$ws.Range("J20").Value = 3279
$ws.Range("K20").Value = 4788.8887
$ws.Range("L20").Value = 3279
$ws.Range("M20").Value = -4541.8887
$ws.Range("N20").Value = -3773
$ws.Range("H134").Value = 3768.6428
$ws.Range("I134").Value = 3850.8462
$ws.Range("J134").Value = 2700
$ws.Range("K134").Value = 11552.5386
$ws.Range("L134").Value = 8100
$ws.Range("M134").Value = -9017.5386
$ws.Range("N134").Value = -13170

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2414.182
$ws.Range("I31").Value = 2505.7
$ws.Range("J31").Value = 1499
$ws.Range("K31").Value = 2505.7
$ws.Range("L31").Value = 1499
$ws.Range("M31").Value = -2210.7
$ws.Range("H34").Value = 2414.182
$ws.Range("I34").Value = 2505.7
$ws.Range("J34").Value = 1499
$ws.Range("K34").Value = 2505.7
$ws.Range("L34").Value = 1499
$ws.Range("M34").Value = -2303.7
$ws.Range("H58").Value = 2715.2
$ws.Range("I58").Value = 2670.4614
$ws.Range("J58").Value = 3006
$ws.Range("K58").Value = 2670.4614
$ws.Range("L58").Value = 3006
$ws.Range("M58").Value = -2467.4614
$ws.Range("N58").Value = -3412
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("H132").Value = 4407.3335
$ws.Range("I132").Value = 4407.3335
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13222.0005
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10692.0005
$ws.Range("H134").Value = 3199.8
$ws.Range("I134").Value = 3199.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9599.400000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7064.400000000001
$ws.Range("H136").Value = 2715.2
$ws.Range("I136").Value = 2670.4614
$ws.Range("J136").Value = 3006
$ws.Range("K136").Value = 8011.3842
$ws.Range("L136").Value = 9018
$ws.Range("M136").Value = -5461.3842
$ws.Range("N136").Value = -14118

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1357.125
$ws.Range("I33").Value = 137
$ws.Range("J33").Value = 3390.6667
$ws.Range("K33").Value = 822
$ws.Range("L33").Value = 20344.0002
$ws.Range("M33").Value = -539
$ws.Range("N33").Value = -20910.0002
$ws.Range("N123").ClearContents()
$ws.Range("H123").Value = 2030
$ws.Range("I123").Value = 2030
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 6090
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -3640
$ws.Range("H141").Value = 14249.833
$ws.Range("I141").Value = 3499.6667
$ws.Range("J141").Value = 25000
$ws.Range("K141").Value = 10499.0001
$ws.Range("L141").Value = 75000
$ws.Range("M141").Value = -5319.000100000001
$ws.Range("N141").Value = -85360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11766.223
$ws.Range("I70").Value = 10999.429
$ws.Range("J70").Value = 14450
$ws.Range("K70").Value = 10999.429
$ws.Range("L70").Value = 14450
$ws.Range("M70").Value = -10729.429
$ws.Range("H73").Value = 11766.223
$ws.Range("I73").Value = 10999.429
$ws.Range("J73").Value = 14450
$ws.Range("K73").Value = 10999.429
$ws.Range("L73").Value = 14450
$ws.Range("M73").Value = -10063.429
$ws.Range("H80").Value = 5145.5
$ws.Range("I80").Value = 2626
$ws.Range("J80").Value = 7665
$ws.Range("K80").Value = 2626
$ws.Range("L80").Value = 7665
$ws.Range("M80").Value = -1628
$ws.Range("N80").Value = -9661
$ws.Range("H83").Value = 5145.5
$ws.Range("I83").Value = 2626
$ws.Range("J83").Value = 7665
$ws.Range("K83").Value = 13130
$ws.Range("L83").Value = 38325
$ws.Range("M83").Value = -8138
$ws.Range("N83").Value = -48309
$ws.Range("N97").ClearContents()
$ws.Range("H97").Value = 2023.8
$ws.Range("I97").Value = 2023.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2023.8
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1527.8
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 1407.8
$ws.Range("I132").Value = 1407.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4223.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1693.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 951
$ws.Range("I22").Value = 940
$ws.Range("J22").Value = 962
$ws.Range("K22").Value = 940
$ws.Range("L22").Value = 962
$ws.Range("M22").Value = -645
$ws.Range("N22").Value = -1552
$ws.Range("H27").Value = 951
$ws.Range("I27").Value = 940
$ws.Range("J27").Value = 962
$ws.Range("K27").Value = 940
$ws.Range("L27").Value = 962
$ws.Range("M27").Value = -833
$ws.Range("N27").Value = -1176
$ws.Range("H40").Value = 4075.75
$ws.Range("I40").Value = 3768
$ws.Range("J40").Value = 4999
$ws.Range("K40").Value = 3768
$ws.Range("L40").Value = 4999
$ws.Range("M40").Value = -3632
$ws.Range("H55").Value = 224.83333
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 299.33334
$ws.Range("K55").Value = 200
$ws.Range("L55").Value = 299.33334
$ws.Range("M55").Value = -27
$ws.Range("H122").Value = 2876.3125
$ws.Range("I122").Value = 2874.4
$ws.Range("J122").Value = 2905
$ws.Range("K122").Value = 8623.200000000001
$ws.Range("L122").Value = 8715
$ws.Range("M122").Value = -6173.200000000001
$ws.Range("H132").Value = 3369.7144
$ws.Range("I132").Value = 2345.75
$ws.Range("J132").Value = 3779.3
$ws.Range("K132").Value = 7037.25
$ws.Range("L132").Value = 11337.9
$ws.Range("M132").Value = -4507.25
$ws.Range("N132").Value = -16397.9
$ws.Range("H136").Value = 3280.2
$ws.Range("I136").Value = 3173.6956
$ws.Range("J136").Value = 4505
$ws.Range("K136").Value = 9521.086800000001
$ws.Range("L136").Value = 13515
$ws.Range("M136").Value = -6971.086800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 404.5
$ws.Range("I107").Value = 348.2857
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 1044.8571
$ws.Range("L107").Value = 2394
$ws.Range("M107").Value = 875.1428999999998
$ws.Range("H122").Value = 1998.2858
$ws.Range("I122").Value = 2016.3334
$ws.Range("J122").Value = 1890
$ws.Range("K122").Value = 6049.0002
$ws.Range("L122").Value = 5670
$ws.Range("M122").Value = -3599.0002
$ws.Range("H132").Value = 2273.2222
$ws.Range("I132").Value = 1738.8422
$ws.Range("J132").Value = 3542.375
$ws.Range("K132").Value = 5216.5266
$ws.Range("L132").Value = 10627.125
$ws.Range("M132").Value = -2686.5266
$ws.Range("N132").Value = -15687.125
$ws.Range("H136").Value = 1247.5625
$ws.Range("I136").Value = 1182.6
$ws.Range("J136").Value = 2222
$ws.Range("K136").Value = 3547.8
$ws.Range("L136").Value = 6666
$ws.Range("M136").Value = -997.7999999999997
$ws.Range("N136").Value = -11766
